# do data cleaning before applying filters
#
# 1. Fix the header text in A1 (drop the stray "번호" suffix).
# 2. Convert the budget column (I) from comma-formatted text strings
#    into real numeric values, for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header cleanup -------------------------------------------------
$ws.Range("A1").Value = "용역 발주계획목록"

# --- 2. Clean up column I (예산액(원)) -----------------------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $clean = $text -replace ",", ""
    $clean = $clean.Trim()

    if ($clean -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = [double]$clean
    }
}
